# ContactData.xlsx edit: add an "Edit Data" sheet that holds the raw
# profile fields (now updated for Riyas Moosa / Designer), and turn the
# ContactData sheet's data row into formulas that pull from it.

$wb = $excel.ActiveWorkbook
$contact = $wb.Worksheets.Item("ContactData")

# ---------------------------------------------------------------------
# 1. Insert the new "Edit Data" worksheet right after "ContactData".
# ---------------------------------------------------------------------
$editData = $wb.Worksheets.Add([Type]::Missing, $contact)
$editData.Name = "Edit Data"

# ---------------------------------------------------------------------
# 2. Header row (row 1) on "Edit Data".
# ---------------------------------------------------------------------
$editData.Range("A1").Value = "name"
$editData.Range("B1").Value = "title"
$editData.Range("C1").Value = "mainImageURI"
$editData.Range("D1").Value = "whatsappNumber"
$editData.Range("E1").Value = "emailAddress"
$editData.Range("F1").Value = "qrCodeUrl"
$editData.Range("G1").Value = "linkedin"
$editData.Range("H1").Value = "instagram"
$editData.Range("I1").Value = "behance"
$editData.Range("J1").Value = "youtube"
$editData.Range("K1").Value = "facebook"
$editData.Range("L1").Value = "webLink1_text"
$editData.Range("M1").Value = "Header BG"
$editData.Range("N1").Value = "WeChat QR"
$editData.Range("O1").Value = "Brand logos"
$editData.Range("P1").Value = "Bio"

# ---------------------------------------------------------------------
# 3. Data row (row 2) on "Edit Data".
# ---------------------------------------------------------------------
$editData.Range("A2").Value = "Riyas Moosa"
$editData.Range("B2").Value = "Designer"
$editData.Range("C2").Value = "https://rmoosa2014.github.io/Resume/DP.png"
$editData.Range("D2").Value = " 966 535531913"
$editData.Range("E2").Value = "ahmed.maher@example.com"
$editData.Range("F2").Value = "https://rmoosa2014.github.io/Resume/"
$editData.Range("G2").Value = "https://sa.linkedin.com/in/riyas-moosa-72923992?trk=people-guest_people_search-card"
$editData.Range("H2").Value = "https://instagram.com/riyas"
$editData.Range("I2").Value = "https://behance.net/riyas"
$editData.Range("J2").Value = "https://youtube.com/@riyas"
$editData.Range("K2").Value = "https://facebook.com/riyas``"
$editData.Range("L2").Value = "www.lightingstores.com"
$editData.Range("M2").Value = "https://i.imgur.com/otJ9G8X.png"
$editData.Range("N2").Value = "https://rmoosa2014.github.io/Resume/"
$editData.Range("O2").Value = "https://rmoosa2014.github.io/Resume/LS LOGO.png"
$editData.Range("P2").Value = "As a versatile Digital Content Creator, I bring ideas to life across the digital landscape, specializing in end-to-end media production, from compelling Design and professional-grade Video Editing to pioneering the use of AI content generation. I thrive on transforming complex concepts into captivating, high-impact digital experiences that drive engagement and tell unforgettable stories."

# Extra webLink2 rows (L3:L5) and extra brand-logo rows (O3:O5).
$editData.Range("L3").Value = "https://lightingstores.com.sa/en"
$editData.Range("L4").Value = "www.iluslighting.com"
$editData.Range("L5").Value = "https://lightingstores.com.sa/en"
$editData.Range("O3").Value = "https://rmoosa2014.github.io/Resume/Illuslogo.svg"
$editData.Range("O4").Value = "https://rmoosa2014.github.io/Resume/Illictlogo.svg"
$editData.Range("O5").Value = "https://rmoosa2014.github.io/Resume/HYP_Logo.png,"

# Spacer cell seen in the original file (Q9, empty but formatted).
$editData.Range("Q9").Value = ""

# ---------------------------------------------------------------------
# 4. Formatting on "Edit Data".
# ---------------------------------------------------------------------
# Row 2 (the data row) is a tall, wrapped/vertically centered row.
$editData.Range("A2:G2").VerticalAlignment = -4108
$editData.Range("L2:N2").VerticalAlignment = -4108
$editData.Range("H2:K2").VerticalAlignment = -4108
$editData.Range("O2").VerticalAlignment = -4108
$editData.Range("P2").WrapText = $true
$editData.Rows.Item(2).RowHeight = 86.4

# ---------------------------------------------------------------------
# 5. Hyperlinks on "Edit Data" (social links + brand logo images).
# ---------------------------------------------------------------------
$editData.Hyperlinks.Add($editData.Range("H2"), "https://instagram.com/riyas") | Out-Null
$editData.Hyperlinks.Add($editData.Range("I2"), "https://behance.net/riyas") | Out-Null
$editData.Hyperlinks.Add($editData.Range("J2"), "https://youtube.com/@riyas") | Out-Null
$editData.Hyperlinks.Add($editData.Range("K2"), "https://facebook.com/riyas``") | Out-Null
$editData.Hyperlinks.Add($editData.Range("O2"), "https://rmoosa2014.github.io/Resume/LS LOGO.png") | Out-Null
$editData.Hyperlinks.Add($editData.Range("O3"), "https://rmoosa2014.github.io/Resume/Illuslogo.svg") | Out-Null
$editData.Hyperlinks.Add($editData.Range("O4"), "https://rmoosa2014.github.io/Resume/Illictlogo.svg") | Out-Null
$editData.Hyperlinks.Add($editData.Range("O5"), "https://rmoosa2014.github.io/Resume/HYP_Logo.png,") | Out-Null

# Re-apply the vertical-center alignment to the hyperlinked cells (adding
# a hyperlink resets the font/style of a cell).
$editData.Range("H2:K2").VerticalAlignment = -4108
$editData.Range("O2").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# 6. Column widths / view on "Edit Data".
# ---------------------------------------------------------------------
$editData.Columns.Item(1).ColumnWidth = 10.26
$editData.Columns.Item(3).ColumnWidth = 17.38
$editData.Columns.Item(4).ColumnWidth = 14.06
$editData.Columns.Item(5).ColumnWidth = 24.48
$editData.Columns.Item(6).ColumnWidth = 32.26
$editData.Columns.Item(7).ColumnWidth = 12.11
$editData.Columns.Item(8).ColumnWidth = 25.16
$editData.Columns.Item(9).ColumnWidth = 22.79
$editData.Columns.Item(10).ColumnWidth = 24.9
$editData.Columns.Item(11).ColumnWidth = 24.38
$editData.Columns.Item(12).ColumnWidth = 26.85
$editData.Columns.Item(13).ColumnWidth = 27
$editData.Columns.Item(14).ColumnWidth = 32.26
$editData.Columns.Item(15).ColumnWidth = 44.48
$editData.Columns.Item(16).ColumnWidth = 57.13

$editData.Range("M2").Select()
$editData.Application.ActiveWindow.ScrollColumn = 6

# ---------------------------------------------------------------------
# 7. ContactData!row1 header cells keep their text (shared-string ids
#    shuffle automatically -- only the visible text matters).
# ---------------------------------------------------------------------
$contact.Range("A1").Value = "name"
$contact.Range("B1").Value = "title"
$contact.Range("C1").Value = "mainImageURI"
$contact.Range("D1").Value = "whatsappNumber"
$contact.Range("E1").Value = "emailAddress"
$contact.Range("F1").Value = "qrCodeUrl"
$contact.Range("G1").Value = "linkedin"
$contact.Range("H1").Value = "instagram"
$contact.Range("I1").Value = "behance"
$contact.Range("J1").Value = "youtube"
$contact.Range("K1").Value = "facebook"
$contact.Range("L1").Value = "webLink1_text"
$contact.Range("M1").Value = "webLink1_href"
$contact.Range("N1").Value = "webLink2_text"
$contact.Range("O1").Value = "webLink2_href"
$contact.Range("P1").Value = "Header BG"
$contact.Range("Q1").Value = "WeChat QR"
$contact.Range("R1").Value = "Brand logos"
$contact.Range("S1").Value = "Bio"

# ---------------------------------------------------------------------
# 8. ContactData!row2 becomes formulas referencing "Edit Data", and its
#    old per-cell hyperlinks + fill style go away.
# ---------------------------------------------------------------------
$contact.Range("A2:S2").Style = "Normal"
$contact.Hyperlinks.Delete()

$contact.Range("A2").Formula = "='Edit Data'!A2"
$contact.Range("B2").Formula = "='Edit Data'!B2"
$contact.Range("C2").Formula = "='Edit Data'!C2"
$contact.Range("D2").Formula = "='Edit Data'!D2"
$contact.Range("E2").Formula = "='Edit Data'!E2"
$contact.Range("F2").Formula = "='Edit Data'!F2"
$contact.Range("G2").Formula = "='Edit Data'!G2"
$contact.Range("H2").Formula = "='Edit Data'!H2"
$contact.Range("I2").Formula = "='Edit Data'!I2"
$contact.Range("J2").Formula = "='Edit Data'!J2"
$contact.Range("K2").Formula = "='Edit Data'!K2"
$contact.Range("L2").Formula = "='Edit Data'!L2"
$contact.Range("M2").Formula = "='Edit Data'!L3"
$contact.Range("N2").Formula = "='Edit Data'!L4"
$contact.Range("O2").Formula = "='Edit Data'!L5"
$contact.Range("P2").Formula = "='Edit Data'!M2"
$contact.Range("Q2").Formula = "='Edit Data'!N2"
$contact.Range("R2").Formula = "='Edit Data'!O2&`",`"&'Edit Data'!O3&`",`"&'Edit Data'!O4&`",`"&'Edit Data'!O5&`",`"&'Edit Data'!O6"
$contact.Range("S2").Formula = "='Edit Data'!P2"

# ---------------------------------------------------------------------
# 9. ContactData column widths + view.
# ---------------------------------------------------------------------
$contact.Columns.Item(12).ColumnWidth = 18.67
$contact.Columns.Item(13).ColumnWidth = 25.17
$contact.Columns.Item(18).ColumnWidth = 168.5
$contact.Columns.Item(19).ColumnWidth = 46.33

$contact.Activate()
$contact.Range("R16").Select()
$contact.Application.ActiveWindow.ScrollColumn = 13
